$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsBVTStL = $wb.Worksheets.Item("BVTStL")

# --- About sheet ---
# Replace the 5-line note (rows 15-19) with a single consolidated line in row 15,
# then clear the now-unused rows 16-19.
$wsAbout.Range("A15").Value = "Based on the California LCFS, we choose to exempt aircraft."
$wsAbout.Range("A16:B19").Clear()

# --- BVTStL sheet ---
# Add a "(Boolean)" label in A1
$wsBVTStL.Range("A1").Value = "(Boolean)"

# Update rail (row 5) and ships (row 6) to be subject to LCFS (1,1) instead of (0,0)
$wsBVTStL.Range("B5").Value = 1
$wsBVTStL.Range("C5").Value = 1
$wsBVTStL.Range("B6").Value = 1
$wsBVTStL.Range("C6").Value = 1

$wb.Save()
